# daily auto push: 2026-02-19 03:14 UTC
# Insert a new data row for 2026/02/19 10:00 (weekday 木) right after the
# existing 2026/02/19 00:00 row (row 830), pushing every following row
# down by one. Excel's row-insert takes care of re-indexing the rest of
# the sheet (and the printed <dimension> ref) automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 830

# Push row 830..871 down to 831..872, leaving a blank row 830 behind.
$ws.Rows.Item($newRow).Insert()

# Column A holds date-like text (e.g. "2026/02/19") that must stay a
# literal string, not get auto-converted into a date serial number by
# Excel's smart-entry parser. Force text format while assigning it, then
# drop the formatting override again so the cell ends up styled exactly
# like its neighbours (no explicit style index).
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "2026/02/19"
$ws.Range("A$newRow").ClearFormats()

$ws.Range("B$newRow").Value = "木"
$ws.Range("C$newRow").Value = 10
$ws.Range("D$newRow").Value = 201
